$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 361, shifting the existing rows 361-378 down to 362-379.
$ws.Rows.Item(361).Insert()

# Populate the newly inserted row 361 with the same static/reference data as the
# (now shifted-down) row that used to occupy 361, but with the new measurement
# values from this week's entry.
$ws.Range("A361").Value = 3
$ws.Range("B361").Value = "Femacal de La Calera"
$ws.Range("C361").Value = "Coquimbo"
$ws.Range("D361").Value = 44753
$ws.Range("E361").Value = 5
$ws.Range("F361").Value = 100112043
$ws.Range("G361").Value = "Pepino ensalada"
$ws.Range("H361").Value = "Sin especificar"
$ws.Range("I361").Value = "Primera"
$ws.Range("J361").Value = 105
$ws.Range("K361").Value = 16000
$ws.Range("L361").Value = 17000
$ws.Range("M361").Value = 16476
$ws.Range("N361").Value = "$/caja 70 unidades"
$ws.Range("O361").Value = "Región de Arica y Parinacota"
$ws.Range("P361").Value = 235
$ws.Range("Q361").Value = 70
$ws.Range("R361").Value = "Hortaliza"
